$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AKS")

# Insert a new column before the existing "Service IP" column (column E)
# so the new "Service Type" column lands at E and "Service IP" shifts to F.
$ws.Columns("E").Insert()

# Header
$ws.Range("E1").Value = "Service Type"

# Fill every data row (2-11) with "ClusterIP" for the new Service Type column
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).Value = "ClusterIP"
}

$ws.Range("E1").ColumnWidth = 10.3

$ws.Range("G:G").Select()
